$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The daily report table repeats a fixed 9-row task block for each date,
# starting at row 2 (row 1 is the header). The dates are being rolled
# forward by one week: 2023-09-18..2023-09-23 -> 2023-09-25..2023-09-30.
$dates = @("2023-09-25", "2023-09-26", "2023-09-27", "2023-09-28", "2023-09-29", "2023-09-30")
$blockSize = 9
$startRow = 2

for ($i = 0; $i -lt $dates.Length; $i++) {
    $rowStart = $startRow + ($i * $blockSize)
    $rowEnd = $rowStart + $blockSize - 1
    $colA = $ws.Cells.Item($rowStart, 1)
    $colAEnd = $ws.Cells.Item($rowEnd, 1)
    $ws.Range($colA, $colAEnd).Value = $dates[$i]
}

# Move the active selection from B47 to B44.
$ws.Range("B44").Select()
